$wb = $excel.ActiveWorkbook

$wsOriginal = $wb.Worksheets.Item("original")
$wsBoxes = $wb.Worksheets.Item("boxes")

# Replace the en-dash ("\u2013") with a plain hyphen in the milestone captions
# on the "original" sheet. The "boxes" sheet pulls these values in via
# formulas that reference "original", so they will recalculate automatically.
$wsOriginal.Range("F3").Value = "Distance schooling for 7-19 year-olds, nationwide"
$wsOriginal.Range("F4").Value = "Distance schooling for 13-19 year-olds in various regions"
$wsOriginal.Range("F5").Value = "Distance schooling for 16-19 year-olds, Helsinki region"
$wsOriginal.Range("F6").Value = "Distance schooling for 13-19 year-olds, Helsinki and certain other regions"
$wsOriginal.Range("F7").Value = "Distance schooling for 13-19 year-olds, nationwide"

$excel.Calculate()

# Make "original" the active sheet with F8 selected, leaving "boxes" with
# F3 last-selected but not the active tab.
$wsBoxes.Activate()
$wsBoxes.Range("F3").Select() | Out-Null

$wsOriginal.Activate()
$wsOriginal.Range("F8").Select() | Out-Null
